# Applies the workbook edit:
#  - Summary sheet: remove the leftover/junk rows 8-13 below the main table
#  - Repayment Schedule sheet: shift the last ("Over Due") column from O to P
#  - Transactions sheet: remove leftover/junk rows (5, 6, 21) and stray cells
#    in columns K:M that sat past the real A1:J3 table
#  - Refresh each affected sheet's selection to match the new layout

$wb = $excel.ActiveWorkbook

# ---- Summary sheet --------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Rows("8:13").Clear()
$wsSummary.Range("C4").Select()

# ---- Repayment Schedule sheet ---------------------------------------
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")
$wsSchedule.Range("O1:O14").Cut($wsSchedule.Range("P1:P14"))
$wsSchedule.Range("G11").Select()

# ---- Transactions sheet ----------------------------------------------
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Range("K2:M2").Clear()
$wsTransactions.Range("K3:M3").Clear()
$wsTransactions.Rows("5:6").Clear()
$wsTransactions.Rows("21:21").Clear()
